$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value2
    $cVal = $cCell.Value2

    if ([string]::IsNullOrEmpty($bVal)) { continue }

    $prefix = $null
    if ($bVal -eq "club-sports") { $prefix = "club" }
    elseif ($bVal -eq "uil-sports") { $prefix = "uil" }

    if ($prefix -ne $null) {
        $suffix = $null
        if ($cVal -eq "Soccer-Boys") { $suffix = "boys" }
        elseif ($cVal -eq "Soccer-Girls") { $suffix = "girls" }
        elseif ($cVal -eq "Soccer-Coed") { $suffix = "coed" }

        if ($suffix -ne $null) {
            $bCell.Value2 = "sports_" + $prefix + "_" + $suffix
        }
    }

    if ($cVal -eq "Soccer-Boys" -or $cVal -eq "Soccer-Girls" -or $cVal -eq "Soccer-Coed") {
        $cCell.Value2 = "Cheer/Drill"
    }
}
